# Update the compilation-file citation DOI column (AP) for rows 10-45 so
# that each cell references the cell immediately above it instead of
# holding its own static string value. AP10 now points at AP9, and
# AP11:AP45 each point at the cell directly above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AP10").Formula = "=AP9"
$ws.Range("AP11:AP45").FormulaR1C1 = "=R[-1]C"

# The old per-cell hyperlink that covered AP10:AP45 no longer applies
# now that those cells are formulas, so remove it. Only the hyperlink on
# AP9 (the true source citation) should remain.
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $hl = $ws.Hyperlinks.Item($i)
    $addr = $hl.Range.Address()
    if ($addr -ne '$AP$9') {
        $hl.Delete() | Out-Null
    }
}

# Reflect the author's final selection on the sheet.
$ws.Range("AP10:AP45").Select() | Out-Null
